$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '36.777.31'
$ws.Range('E2').Value = '  +0.91%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.965.95'
$ws.Range('E3').Value = '  +1.14%  '
$ws.Range('E4').Value = '  +0.00%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '244.87'
$ws.Range('E5').Value = '  +0.43%  '
$ws.Range('E6').Value = '  +0.91%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '59.19'
$ws.Range('E7').Value = '  +0.86%  '
$ws.Range('E8').Value = '  -0.03%  '
$ws.Range('E9').Value = '  +1.98%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.0811'
$ws.Range('E10').Value = '  -3.62%  '
$ws.Range('E11').Value = '  -0.34%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '22.53'
$ws.Range('E12').Value = '  +4.43%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '2.250.72'
$ws.Range('E13').Value = '  +1.06%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '0.827'
$ws.Range('E14').Value = '  -0.08%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '13.76'
$ws.Range('E15').Value = '  +0.82%  '
$ws.Range('E16').Value = '  +0.29%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '1.957.91'
$ws.Range('E17').Value = '  -0.48%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '36.675.05'
$ws.Range('E18').Value = '  +0.79%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '69.83'
$ws.Range('E19').Value = '  +0.07%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '0.0₃0864'
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '5.10'
$ws.Range('E21').Value = '  +1.76%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '229.28'
$ws.Range('E22').Value = '  -0.33%  '
$ws.Range('E23').Value = '  -0.09%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '2.45'
$ws.Range('E24').Value = '  -0.31%  '
$ws.Range('E25').Value = '  +2.74%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '9.34'
$ws.Range('E26').Value = '  +0.23%  '
$ws.Range('E27').Value = '  +14.14%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '160.72'
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '19.43'
$ws.Range('E29').Value = '  -0.19%  '
$ws.Range('E30').Value = '  +0.87%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '1.14'
$ws.Range('E31').Value = '  -1.61%  '
$ws.Range('E32').Value = '  +0.62%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.0620'
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '4.29'
$ws.Range('E34').Value = '  +0.06%  '
$ws.Range('E35').Value = '  +0.01%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '6.11'
$ws.Range('E36').Value = '  -2.66%  '
$ws.Range('E37').Value = '  +4.70%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '3.41'
$ws.Range('E38').Value = '  +11.46%  '
$ws.Range('E39').Value = '  -0.14%  '
$ws.Range('E40').Value = '  +3.11%  '
$ws.Range('E41').Value = '  -2.39%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.0213'
$ws.Range('E42').Value = '  +1.67%  '
$ws.Range('E43').Value = '  -1.21%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '16.14'
$ws.Range('E44').Value = '  +0.20%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '1.361.54'
$ws.Range('E45').Value = '  +0.58%  '
$ws.Range('E46').Value = '  +0.38%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '87.90'
$ws.Range('E47').Value = '  +0.04%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '7.17'
$ws.Range('E48').Value = '  -0.51%  '
$ws.Range('E49').Value = '  +0.82%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '2.142.09'
$ws.Range('E50').Value = '  +1.05%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '43.77'
$ws.Range('E51').Value = '  -4.07%  '
